$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("S4").Value = 2022
$ws.Range("S5").Value = 27.292394741221504
